# Populate B1 and A2 with a numeric 0, and B2 with the label
# "disconnected_elements" (via the shared-string table), then apply the
# "header cell" look (bold, centered/top-aligned, thin box border) to the
# two numeric label cells B1 and A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$rng = $ws.Range("B1,A2")
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.VerticalAlignment = -4160     # xlTop
$rng.Borders.LineStyle = 1         # xlContinuous
$rng.Borders.Weight = 2            # xlThin
